# Commit message: "Use another category "HTC""
#
# The sheet "Перечень настроек" contained a placeholder row (row 33) with
# the text "Категории добавить описание" ("Add category description") that
# was meant to be replaced by a note about using category "HTC" as the test
# fixture. That note already existed one row below, in column C of row 34:
# "все эти товары можно взять из категории HTC и работать только с этой
# категорией" ("all these products can be taken from category HTC and work
# only with this category").
#
# The fix deletes the now-obsolete placeholder row. Everything below shifts
# up by one row, so the HTC note ends up sharing a row with "В категории
# должен быть товар без наличия" instead of sitting alone, the sheet
# shrinks from 40 to 39 rows, and focus moves from the "Коды значений"
# sheet to this one, which is also renamed to reflect its broader scope.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename the first sheet to reflect its new, broader purpose.
$ws1.Name = "Настройки и Условия"

# Remove the obsolete placeholder row; everything below it (including the
# HTC note) shifts up by one row and the now-unused shared string is
# dropped automatically.
$ws1.Rows.Item(33).Delete()

# Sheet 1 becomes the active/selected sheet (previously it was sheet 2).
$ws1.Activate()

# Restore a sensible cursor/viewport position on sheet 1.
$ws1.Range("C15").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
